$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell carrying the default (unstyled) format, used to restore
# cell style after forcing numeric-looking text values to stay as text.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '55.911.90'
$ws.Range("E2").Value = '  -0.98%  '

$ws.Range("D3").Value = '2.292.83'
$ws.Range("E3").Value = '  -0.96%  '

$ws.Range("E4").Value = '  -0.02%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '512.58'
$cell.Style = $defaultStyle
$ws.Range("E5").Value = '  -0.95%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '129.82'
$cell.Style = $defaultStyle
$ws.Range("E6").Value = '  -4.02%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $defaultStyle
$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("E8").Value = '  -1.52%  '

$ws.Range("D9").Value = '2.295.99'
$ws.Range("E9").Value = '  -1.40%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0989'
$cell.Style = $defaultStyle
$ws.Range("E10").Value = '  -3.45%  '

$ws.Range("E11").Value = '  -0.16%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '5.23'
$cell.Style = $defaultStyle
$ws.Range("E12").Value = '  -1.79%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.333'
$cell.Style = $defaultStyle
$ws.Range("E13").Value = '  -2.56%  '

$ws.Range("D14").Value = '2.699.80'
$ws.Range("E14").Value = '  -0.89%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '23.15'
$cell.Style = $defaultStyle
$ws.Range("E15").Value = '  -3.39%  '

$ws.Range("D16").Value = '55.860.37'
$ws.Range("E16").Value = '  -1.22%  '

$ws.Range("E17").Value = '  -2.67%  '

$ws.Range("D18").Value = '2.310.51'
$ws.Range("E18").Value = '  -1.08%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '10.30'
$cell.Style = $defaultStyle
$ws.Range("E19").Value = '  -1.83%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '326.40'
$cell.Style = $defaultStyle
$ws.Range("E20").Value = '  +1.15%  '

$ws.Range("E21").Value = '  -2.90%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.63'
$cell.Style = $defaultStyle
$ws.Range("E22").Value = '  +0.57%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $defaultStyle

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '60.36'
$cell.Style = $defaultStyle
$ws.Range("E24").Value = '  -0.17%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.162'
$cell.Style = $defaultStyle
$ws.Range("E25").Value = '  -0.60%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '8.51'
$cell.Style = $defaultStyle
$ws.Range("E26").Value = '  +6.55%  '

$ws.Range("E27").Value = '  +0.52%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.30'
$cell.Style = $defaultStyle
$ws.Range("E28").Value = '  +0.92%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '168.03'
$cell.Style = $defaultStyle
$ws.Range("E29").Value = '  +0.79%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.67'
$cell.Style = $defaultStyle
$ws.Range("E30").Value = '  -1.74%  '

$ws.Range("D31").Value = '0.0₃0705'
$ws.Range("E31").Value = '  -4.69%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '6.04'
$cell.Style = $defaultStyle
$ws.Range("E32").Value = '  -2.67%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '18.14'
$cell.Style = $defaultStyle
$ws.Range("E33").Value = '  -1.23%  '

$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("E35").Value = '  +0.50%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.22'
$cell.Style = $defaultStyle
$ws.Range("E36").Value = '  -2.59%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '3.88'
$cell.Style = $defaultStyle
$ws.Range("E37").Value = '  -3.60%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.877'
$cell.Style = $defaultStyle
$ws.Range("E38").Value = '  -5.08%  '

$ws.Range("E39").Value = '  +0.07%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '38.24'
$cell.Style = $defaultStyle
$ws.Range("E40").Value = '  +0.78%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '146.20'
$cell.Style = $defaultStyle
$ws.Range("E41").Value = '  +4.16%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.370'
$cell.Style = $defaultStyle
$ws.Range("E42").Value = '  -3.05%  '

$ws.Range("E43").Value = '  -1.92%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '279.43'
$cell.Style = $defaultStyle
$ws.Range("E44").Value = '  +0.38%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '4.88'
$cell.Style = $defaultStyle
$ws.Range("E45").Value = '  -6.47%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0920'
$cell.Style = $defaultStyle
$ws.Range("E46").Value = '  -1.30%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.0491'
$cell.Style = $defaultStyle
$ws.Range("E47").Value = '  -3.38%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.551'
$cell.Style = $defaultStyle
$ws.Range("E48").Value = '  -1.90%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '17.95'
$cell.Style = $defaultStyle
$ws.Range("E49").Value = '  +0.71%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0212'
$cell.Style = $defaultStyle
$ws.Range("E50").Value = '  -2.61%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '16.85'
$cell.Style = $defaultStyle
$ws.Range("E51").Value = '  +0.21%  '

